$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.646.92"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "3.042.18"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "555.95"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.040.63"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "0.519"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  -13.09%  "
$ws.Range("D12").Value = "0.484"
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "35.40"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "3.538.85"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "63.701.85"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "3.039.15"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "479.09"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "14.52"
$ws.Range("E23").Value = "  +8.45%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "82.48"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "8.04"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "26.05"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").Value = "54.93"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "439.70"
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("D39").Value = "0.0812"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").Value = "2.994.25"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").Value = "8.28"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "0.269"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").Value = "27.70"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "2.22"
$ws.Range("E46").Value = "  +5.17%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "117.67"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "2.08"
$ws.Range("E51").Value = "  +0.10%  "
